# Apply updated Betfair back/lay odds values to Sheet1 (rows 2-12).
# Values below correspond exactly to the cells changed in the source commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.26
$ws.Range("G2").Value = 2.4
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 4.6
$ws.Range("K2").Value = 3.25
$ws.Range("N2").Value = 2.36
$ws.Range("O2").Value = 1.64
$ws.Range("P2").Value = 1.44
$ws.Range("Q2").Value = 2.84
$ws.Range("T2").Value = 2.28
$ws.Range("U2").Value = 1.64
$ws.Range("V2").Value = 1.28
$ws.Range("W2").Value = 1.71
$ws.Range("AC2").Value = 14
$ws.Range("AD2").Value = 980
$ws.Range("AF2").Value = 1000
$ws.Range("AH2").Value = 1000

# Row 3
$ws.Range("F3").Value = 2.18
$ws.Range("G3").Value = 2.28
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3.75
$ws.Range("J3").Value = 3.35
$ws.Range("K3").Value = 3.7
$ws.Range("L3").Value = 1.65
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 2.36
$ws.Range("O3").Value = 1.61
$ws.Range("P3").Value = 1.49
$ws.Range("Q3").Value = 2.9
$ws.Range("S3").Value = 6.4
$ws.Range("T3").Value = 2.36
$ws.Range("U3").Value = 1.59
$ws.Range("V3").Value = 1.37
$ws.Range("W3").Value = 1.71
$ws.Range("X3").Value = 9
$ws.Range("Y3").Value = 8.8
$ws.Range("Z3").Value = 65
$ws.Range("AA3").Value = 100
$ws.Range("AB3").Value = 8
$ws.Range("AC3").Value = 9.8
$ws.Range("AD3").Value = 16.5
$ws.Range("AE3").Value = 75
$ws.Range("AF3").Value = 12.5
$ws.Range("AG3").Value = 14
$ws.Range("AH3").Value = 980
$ws.Range("AI3").Value = 110
$ws.Range("AJ3").Value = 36
$ws.Range("AK3").Value = 44
$ws.Range("AL3").Value = 95
$ws.Range("AN3").Value = 980
$ws.Range("AO3").Value = 120

# Row 4
$ws.Range("G4").Value = 1.17
$ws.Range("H4").Value = 21
$ws.Range("I4").Value = 25
$ws.Range("K4").Value = 11.5
$ws.Range("N4").Value = 10
$ws.Range("Q4").Value = 1.29
$ws.Range("R4").Value = 2.16
$ws.Range("S4").Value = 1.76
$ws.Range("T4").Value = 1.94
$ws.Range("U4").Value = 1.9
$ws.Range("X4").Value = 80
$ws.Range("Z4").Value = 280
$ws.Range("AB4").Value = 17.5
$ws.Range("AC4").Value = 40
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 370
$ws.Range("AJ4").Value = 9.8
$ws.Range("AL4").Value = 36
$ws.Range("AN4").Value = 2.66

# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("K5").Value = 3.25
$ws.Range("L5").Value = 1.68
$ws.Range("N5").Value = 2.24
$ws.Range("Q5").Value = 3.05
$ws.Range("R5").Value = 1.13
$ws.Range("U5").Value = 1.54
$ws.Range("Z5").Value = 44
$ws.Range("AD5").Value = 28
$ws.Range("AH5").Value = 36
$ws.Range("AK5").Value = 36
$ws.Range("AL5").Value = 160

# Row 6
$ws.Range("F6").Value = 2.62
$ws.Range("G6").Value = 2.64
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 3.45
$ws.Range("J6").Value = 3.05
$ws.Range("K6").Value = 3.1
$ws.Range("N6").Value = 2.82
$ws.Range("P6").Value = 1.6
$ws.Range("Q6").Value = 2.64
$ws.Range("V6").Value = 1.4
$ws.Range("W6").Value = 1.6
$ws.Range("X6").Value = 8.8

# Row 7
$ws.Range("F7").Value = 1.61
$ws.Range("G7").Value = 1.62
$ws.Range("N7").Value = 7.8
$ws.Range("O7").Value = 1.12
$ws.Range("P7").Value = 3.3
$ws.Range("Q7").Value = 1.38
$ws.Range("S7").Value = 1.99
$ws.Range("T7").Value = 1.53
$ws.Range("U7").Value = 2.68
$ws.Range("V7").Value = 1.21
$ws.Range("W7").Value = 2.6
$ws.Range("X7").Value = 46
$ws.Range("AA7").Value = 130
$ws.Range("AB7").Value = 16.5
$ws.Range("AF7").Value = 15
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 17
$ws.Range("AJ7").Value = 18
$ws.Range("AN7").Value = 5.1
$ws.Range("AO7").Value = 970

# Row 8
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 5.2
$ws.Range("H8").Value = 1.85
$ws.Range("I8").Value = 1.9
$ws.Range("K8").Value = 3.85
$ws.Range("N8").Value = 3.7
$ws.Range("P8").Value = 1.9
$ws.Range("R8").Value = 1.34
$ws.Range("S8").Value = 3.4
$ws.Range("T8").Value = 1.84
$ws.Range("V8").Value = 2.1
$ws.Range("W8").Value = 1.23
$ws.Range("X8").Value = 15.5
$ws.Range("Y8").Value = 9.2
$ws.Range("AB8").Value = 27
$ws.Range("AC8").Value = 9.2
$ws.Range("AE8").Value = 21
$ws.Range("AF8").Value = 95
$ws.Range("AI8").Value = 95
$ws.Range("AK8").Value = 320
$ws.Range("AL8").Value = 330

# Row 9
$ws.Range("F9").Value = 1.7
$ws.Range("G9").Value = 1.78
$ws.Range("K9").Value = 3.7
$ws.Range("M9").Value = 1.12
$ws.Range("N9").Value = 2.56
$ws.Range("O9").Value = 1.53
$ws.Range("P9").Value = 1.53
$ws.Range("Q9").Value = 2.56
$ws.Range("T9").Value = 2.34
$ws.Range("U9").Value = 1.62
$ws.Range("W9").Value = 2.26
$ws.Range("Y9").Value = 16.5
$ws.Range("Z9").Value = 60
$ws.Range("AJ9").Value = 20

# Row 10
$ws.Range("I10").Value = 10.5
$ws.Range("N10").Value = 3.4
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 1.83
$ws.Range("Q10").Value = 2.18
$ws.Range("T10").Value = 2.44
$ws.Range("U10").Value = 1.67
$ws.Range("V10").Value = 1.1
$ws.Range("AA10").Value = 450
$ws.Range("AC10").Value = 10
$ws.Range("AG10").Value = 10.5
$ws.Range("AH10").Value = 32
$ws.Range("AK10").Value = 17.5
$ws.Range("AN10").Value = 9.2
$ws.Range("AO10").Value = 370

# Row 11
$ws.Range("H11").Value = 6
$ws.Range("O11").Value = 1.31
$ws.Range("P11").Value = 2.06
$ws.Range("Q11").Value = 1.91
$ws.Range("T11").Value = 1.95
$ws.Range("U11").Value = 2.04
$ws.Range("Y11").Value = 19.5
$ws.Range("AC11").Value = 9.2
$ws.Range("AE11").Value = 85
$ws.Range("AF11").Value = 9
$ws.Range("AG11").Value = 9.6
$ws.Range("AH11").Value = 21
$ws.Range("AI11").Value = 75
$ws.Range("AN11").Value = 9.2

# Row 12
$ws.Range("F12").Value = 2.14
$ws.Range("G12").Value = 2.22
$ws.Range("H12").Value = 4.1
$ws.Range("I12").Value = 4.4
$ws.Range("K12").Value = 3.4
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 2.98
$ws.Range("P12").Value = 1.65
$ws.Range("S12").Value = 4.5
$ws.Range("T12").Value = 1.99
$ws.Range("U12").Value = 1.87
$ws.Range("V12").Value = 1.29
$ws.Range("AA12").Value = 270
$ws.Range("AB12").Value = 8
$ws.Range("AE12").Value = 65
$ws.Range("AH12").Value = 22
$ws.Range("AI12").Value = 190
$ws.Range("AJ12").Value = 29
$ws.Range("AL12").Value = 55
$ws.Range("AN12").Value = 29
